# Fixed an issue with digit_count
# Replace the dummy A2:A41 sample data with the real dataset (A2:A59),
# apply number formatting / right alignment / thin left-right borders
# (plus a bottom border under the final row) with the cells unlocked,
# and update the sheet selection/view to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New data values for A2:A59 --------------------------------------
$values = @(
    1680380, 1204, 40532, 210135, 45346, 21858, 1165986, 27601, 191220,
    1009503, 28896, 136261, 179595, 19001, 905828, 152790, 68195, 32317,
    9992236, 156343, 261282, 17118, 91361, 281814, 8663, 13226, 438390,
    137485, 102249, 3185516, 405889, 19749, 2422847, 1586465, 64521,
    2183239, 3296045, 870393, 780558, 281879, 762511, 448244, 1931026,
    270462, 181984, 3228, 43991, 452698, 488281, 553217, 99464, 65682,
    16089, 473891, 55374, 843310, 216291, 81994
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $values[$i]
}

# --- Formatting --------------------------------------------------------
# Build the formatting on hidden helper cells first and copy/paste the
# resulting formats onto the real ranges in a single shot each, so that
# every cell in a given range picks up exactly one final style instead
# of accumulating a pile of unused intermediate styles.

$bodyTemplate = $ws.Range("C1")
$bodyTemplate.NumberFormat = "#,##0"
$bodyTemplate.HorizontalAlignment = -4152
$bodyTemplate.Locked = $false
$bodyTemplate.Borders.Item(7).LineStyle = 1
$bodyTemplate.Borders.Item(10).LineStyle = 1

$bodyTemplate.Copy()
$lastTemplate = $ws.Range("C2")
$lastTemplate.PasteSpecial(-4122)
$lastTemplate.Borders.Item(9).LineStyle = 1

$bodyTemplate.Copy()
$ws.Range("A2:A58").PasteSpecial(-4122)

$lastTemplate.Copy()
$ws.Range("A59").PasteSpecial(-4122)

$bodyTemplate.Clear()
$lastTemplate.Clear()

# --- Sheet view / selection --------------------------------------------
$ws.Activate()
$ws.Range("A2:A59").Select()
$excel.ActiveWindow.ScrollRow = 48
